# The document's cover/date line reads "Septiembre 2016" (paragraph
# styled "Date"). The commit bumps this to "Enero 2017" — update the
# month and year in place.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Septiembre", $true, $false, $false, $false, $false, `
               $true, 1, $false, "Enero", 2)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("2016", $true, $false, $false, $false, $false, `
                $true, 1, $false, "2017", 2)
